$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.113.09'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.49'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4587'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.61%  '
$ws.Range("E8").Value = '  +1.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07333'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8605'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.00'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.823.56'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.699'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.353'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07084'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008850'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.130.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.196'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.224'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.65%  '
$ws.Range("E27").Value = '  +1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.268'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08897'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.7734'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.195'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.974'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.473'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.104'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("E37").Value = '  +0.86%  '
$ws.Range("E38").Value = '  +0.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5401'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.78%  '
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.886'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1715'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5270'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +12.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.633'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.989'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +10.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06494'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.62%  '
$ws.Range("E49").Value = '  +1.42%  '
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9260'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.68%  '
